$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the empty column C, shifting the whole board one column to the left.
$ws.Columns("C").Delete()

# The "position index" columns/rows used to count 1..7; the rule now starts
# at position 0, so decrement every index value by one.
for ($r = 4; $r -le 10; $r++) {
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($r, 3).Value2 - 1   # column C
    $ws.Cells.Item($r, 12).Value = $ws.Cells.Item($r, 12).Value2 - 1 # column L
}

for ($c = 4; $c -le 10; $c++) {
    $ws.Cells.Item(11, $c).Value = $ws.Cells.Item(11, $c).Value2 - 1       # D11:J11
    $ws.Cells.Item(11, $c + 9).Value = $ws.Cells.Item(11, $c + 9).Value2 - 1 # M11:S11
}

# Update the selected cell to reflect where the player now starts.
$ws.Range("J7").Select()
